$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Mentioned_in_text")
$ws1.Range("A15").Value = "list of all genes with their respective lfp values"

$ws2 = $wb.Worksheets.Item("Extra_on_github")
$ws2.Range("A6").Value = "also at beginning of each metacell scripts I could add a link to where the raw data can be downloaded.."

[void]$ws2.Range("A7").Select()
[void]$ws1.Range("C15").Select()
